$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 13
$ws.Range("H13").Value = 12943.556
$ws.Range("J13").Value = 11623
$ws.Range("L13").Value = 11623
$ws.Range("N13").Value = -11961
# ALC row 33
$ws.Range("H33").Value = 239.72414
$ws.Range("I33").Value = 239.72414
$ws.Range("K33").Value = 239.72414
$ws.Range("M33").Value = -10.72414000000001
# ALC row 69
$ws.Range("H69").Value = 3619.1667
$ws.Range("J69").Value = 3619.1667
$ws.Range("L69").Value = 10857.5001
$ws.Range("N69").Value = -12605.5001
# ALC row 72
$ws.Range("H72").Value = 3619.1667
$ws.Range("J72").Value = 3619.1667
$ws.Range("L72").Value = 32572.5003
$ws.Range("N72").Value = -41308.5003
# ALC row 100
$ws.Range("H100").Value = 2255.3125
$ws.Range("I100").Value = 1959.875
$ws.Range("J100").Value = 2550.75
$ws.Range("K100").Value = 1959.875
$ws.Range("L100").Value = 2550.75
$ws.Range("M100").Value = -1418.875
$ws.Range("N100").Value = -3632.75
# ALC row 101
$ws.Range("H101").Value = 581.3333
$ws.Range("I101").Value = 693.5
$ws.Range("K101").Value = 2080.5
$ws.Range("M101").Value = -458.5
# ALC row 106
$ws.Range("H106").Value = 1421.5555
$ws.Range("I106").Value = 1018.8
$ws.Range("K106").Value = 1018.8
$ws.Range("M106").Value = -387.8
# ALC row 116
$ws.Range("H116").Value = 100001750
$ws.Range("I116").Value = 200000000
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 200000000
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = -199996558
$ws.Range("N116").Value = -10384
# ALC row 121
$ws.Range("H121").Value = 642.5714
$ws.Range("I121").Value = 1097.5
$ws.Range("J121").Value = 607.5769
$ws.Range("K121").Value = 3292.5
$ws.Range("L121").Value = 1822.7307
$ws.Range("M121").Value = -1545.5
$ws.Range("N121").Value = -5316.7307
# ALC row 125
$ws.Range("H125").Value = 946.0769
$ws.Range("I125").Value = 790.875
$ws.Range("J125").Value = 1194.4
$ws.Range("K125").Value = 7117.875
$ws.Range("L125").Value = 10749.6
$ws.Range("M125").Value = -4657.875
$ws.Range("N125").Value = -15669.6
# ALC row 132
$ws.Range("H132").Value = 1244.6543
$ws.Range("I132").Value = 1070.1487
$ws.Range("K132").Value = 3210.4461
$ws.Range("M132").Value = -680.4461000000001
# ALC row 137
$ws.Range("H137").Value = 1497.9387
$ws.Range("I137").Value = 1286.5405
$ws.Range("K137").Value = 3859.6215
$ws.Range("M137").Value = -1309.6215
# ALC row 140
$ws.Range("H140").Value = 77053.336
$ws.Range("J140").Value = 75553.63
$ws.Range("L140").Value = 75553.63
$ws.Range("N140").Value = -85913.63

$ws = $wb.Worksheets.Item("ARM")
# ARM row 13
$ws.Range("H13").Value = 30000
$ws.Range("I13").Value = 30000
$ws.Range("K13").Value = 30000
$ws.Range("M13").Value = -29856
# ARM row 63
$ws.Range("H63").Value = 11064.286
$ws.Range("J63").Value = 9367.5
$ws.Range("L63").Value = 9367.5
$ws.Range("N63").Value = -10739.5
# ARM row 66
$ws.Range("H66").Value = 11064.286
$ws.Range("J66").Value = 9367.5
$ws.Range("L66").Value = 46837.5
$ws.Range("N66").Value = -53701.5
# ARM row 74
$ws.Range("H74").Value = 1875.2609
$ws.Range("I74").Value = 1913.5
$ws.Range("J74").Value = 1845.8462
$ws.Range("K74").Value = 1913.5
$ws.Range("L74").Value = 1845.8462
$ws.Range("M74").Value = -1039.5
$ws.Range("N74").Value = -3593.8462
# ARM row 77
$ws.Range("H77").Value = 1875.2609
$ws.Range("I77").Value = 1913.5
$ws.Range("J77").Value = 1845.8462
$ws.Range("K77").Value = 9567.5
$ws.Range("L77").Value = 9229.231
$ws.Range("M77").Value = -5199.5
$ws.Range("N77").Value = -17965.231
# ARM row 102
$ws.Range("H102").Value = 1136.8125
$ws.Range("I102").Value = 1034.1428
$ws.Range("J102").Value = 1855.5
$ws.Range("K102").Value = 1034.1428
$ws.Range("L102").Value = 1855.5
$ws.Range("M102").Value = 587.8571999999999
$ws.Range("N102").Value = -5099.5
# ARM row 122
$ws.Range("H122").Value = 2557.25
$ws.Range("I122").Value = 2545.9583
$ws.Range("J122").Value = 2625
$ws.Range("K122").Value = 7637.874899999999
$ws.Range("L122").Value = 7875
$ws.Range("M122").Value = -5187.874899999999
$ws.Range("N122").Value = -12775
# ARM row 132
$ws.Range("H132").Value = 446001.25
$ws.Range("I132").Value = 541611.25
$ws.Range("J132").Value = 3805
$ws.Range("K132").Value = 1624833.75
$ws.Range("L132").Value = 11415
$ws.Range("M132").Value = -1622303.75
$ws.Range("N132").Value = -16475

$ws = $wb.Worksheets.Item("CRP")
# CRP row 11
$ws.Range("H11").Value = 25333.334
$ws.Range("I11").Value = 30000
$ws.Range("J11").Value = 16000
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = -29860
$ws.Range("N11").Value = -16280
# CRP row 31
$ws.Range("H31").Value = 3121.2273
$ws.Range("I31").Value = 1936.3478
$ws.Range("K31").Value = 1936.3478
$ws.Range("M31").Value = -1641.3478
# CRP row 34
$ws.Range("H34").Value = 3121.2273
$ws.Range("I34").Value = 1936.3478
$ws.Range("K34").Value = 1936.3478
$ws.Range("M34").Value = -1734.3478
# CRP row 58
$ws.Range("H58").Value = 1611724.2
$ws.Range("I58").Value = 2316321
$ws.Range("J58").Value = 1217.5714
$ws.Range("K58").Value = 2316321
$ws.Range("L58").Value = 1217.5714
$ws.Range("M58").Value = -2316118
$ws.Range("N58").Value = -1623.5714
# CRP row 62
$ws.Range("H62").Value = 67313.75
$ws.Range("I62").Value = 127377.5
$ws.Range("J62").Value = 7250
$ws.Range("K62").Value = 127377.5
$ws.Range("L62").Value = 7250
$ws.Range("M62").Value = -126753.5
$ws.Range("N62").Value = -8498
# CRP row 65
$ws.Range("H65").Value = 67313.75
$ws.Range("I65").Value = 127377.5
$ws.Range("J65").Value = 7250
$ws.Range("K65").Value = 636887.5
$ws.Range("L65").Value = 36250
$ws.Range("M65").Value = -633767.5
$ws.Range("N65").Value = -42490
# CRP row 134
$ws.Range("H134").Value = 1464.3704
$ws.Range("I134").Value = 1338.6666
$ws.Range("J134").Value = 1904.3334
$ws.Range("K134").Value = 4015.9998
$ws.Range("L134").Value = 5713.0002
$ws.Range("M134").Value = -1480.9998
$ws.Range("N134").Value = -10783.0002
# CRP row 136
$ws.Range("H136").Value = 1611724.2
$ws.Range("I136").Value = 2316321
$ws.Range("J136").Value = 1217.5714
$ws.Range("K136").Value = 6948963
$ws.Range("L136").Value = 3652.7142
$ws.Range("M136").Value = -6946413
$ws.Range("N136").Value = -8752.7142

$ws = $wb.Worksheets.Item("CUL")
# CUL row 129
$ws.Range("H129").Value = 1787498.8
$ws.Range("J129").Value = 2001926.6
$ws.Range("L129").Value = 6005779.800000001
$ws.Range("N129").Value = -6015779.800000001
# CUL row 131
$ws.Range("H131").Value = 12197754
$ws.Range("J131").Value = 14494224
$ws.Range("L131").Value = 43482672
$ws.Range("N131").Value = -43492752

$ws = $wb.Worksheets.Item("GSM")
# GSM row 7
$ws.Range("H7").Value = 451.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 451.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 451.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -675.5
# GSM row 8
$ws.Range("H8").Value = 451.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 451.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 451.5
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -729.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 8472.286
$ws.Range("I40").Value = 8401.333000000001
$ws.Range("K40").Value = 8401.333000000001
$ws.Range("M40").Value = -8265.333000000001
# LTW row 122
$ws.Range("H122").Value = 13976434
$ws.Range("I122").Value = 17469104
$ws.Range("J122").Value = 5751.25
$ws.Range("K122").Value = 52407312
$ws.Range("L122").Value = 17253.75
$ws.Range("M122").Value = -52404862
$ws.Range("N122").Value = -22153.75

$ws = $wb.Worksheets.Item("WVR")
# WVR row 8
$ws.Range("H8").Value = 1000.6667
$ws.Range("J8").Value = 999
$ws.Range("L8").Value = 999
$ws.Range("N8").Value = -1279
# WVR row 107
$ws.Range("H107").Value = 940.2857
$ws.Range("I107").Value = 948
$ws.Range("J107").Value = 930
$ws.Range("K107").Value = 2844
$ws.Range("L107").Value = 2790
$ws.Range("M107").Value = -924
$ws.Range("N107").Value = -6630
# WVR row 123
$ws.Range("H123").Value = 23676.062
$ws.Range("J123").Value = 23676.062
$ws.Range("L123").Value = 23676.062
$ws.Range("N123").Value = -33476.06200000001
# WVR row 132
$ws.Range("H132").Value = 2596.5
$ws.Range("I132").Value = 2360.75
$ws.Range("J132").Value = 2879.4
$ws.Range("K132").Value = 7082.25
$ws.Range("L132").Value = 8638.200000000001
$ws.Range("M132").Value = -4552.25
$ws.Range("N132").Value = -13698.2
# WVR row 136
$ws.Range("H136").Value = 2417.4688
$ws.Range("I136").Value = 2550.6365
$ws.Range("J136").Value = 2124.5
$ws.Range("K136").Value = 7651.9095
$ws.Range("L136").Value = 6373.5
$ws.Range("M136").Value = -5101.9095
$ws.Range("N136").Value = -11473.5
